$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data block for Iteration III burn-down (rows 36-42) ---
$ws.Range("A36").Value = "Iteration II"
$ws.Range("A36").Font.Bold = $true

$ws.Range("B37").Value = "Ideal Tasks Remaining"
$ws.Range("C37").Value = "Actual Tasks Remaining"

$ws.Range("A38").Value = 0.5
$ws.Range("A39").Value = 0.51041666666666663
$ws.Range("A40").Value = 0.52083333333333304
$ws.Range("A41").Value = 0.53125
$ws.Range("A42").Value = 0.54166666666666696
$ws.Range("A38:A42").NumberFormat = "h:mm"

$ws.Range("B38").Value = 4
$ws.Range("B39").Value = 3.1
$ws.Range("B40").Value = 1.8
$ws.Range("B41").Value = 0.9
$ws.Range("B42").Value = 0

$ws.Range("C38").Value = 4
$ws.Range("C39").Value = 4
$ws.Range("C40").Value = 2
$ws.Range("C41").Value = 1
$ws.Range("C42").Value = 0

$ws.Range("B38:C42").NumberFormat = "0.00"

# --- sheet view state ---
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Range("D48").Select()
